$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new date column before DH (col 112) ---
$ws1 = $wb.Worksheets.Item("Prix Spot")

# Inserting the entire column shifts DH:EL -> DI:EM, carrying over styles/values.
$ws1.Range("DH1").EntireColumn.Insert()

# Populate the freshly inserted (now-empty) column with the new date header
# and "-" placeholder values for every data row (2-25), matching the
# existing pattern used for days without data.
$ws1.Range("DH1").Value = "03-nov"
for ($r = 2; $r -le 25; $r++) {
    $ws1.Cells.Item($r, 112).Value = "-"
}

# --- Sheet "Gaz": append two new daily rows ---
$ws2 = $wb.Worksheets.Item("Gaz")

$ws2.Range("A140").NumberFormat = "@"
$ws2.Cells.Item(140, 1).Value = "2025-11-01"
$ws2.Range("A140").Style = "Normal"
$ws2.Cells.Item(140, 2).Value = 29.8

$ws2.Range("A141").NumberFormat = "@"
$ws2.Cells.Item(141, 1).Value = "2025-11-02"
$ws2.Range("A141").Style = "Normal"
$ws2.Cells.Item(141, 2).Value = 29.8

# --- Sheet "CO2": append two new daily rows ---
$ws3 = $wb.Worksheets.Item("CO2")

$ws3.Range("A140").NumberFormat = "@"
$ws3.Cells.Item(140, 1).Value = "2025-11-01"
$ws3.Range("A140").Style = "Normal"
$ws3.Cells.Item(140, 2).Value = 78

$ws3.Range("A141").NumberFormat = "@"
$ws3.Cells.Item(141, 1).Value = "2025-11-02"
$ws3.Range("A141").Style = "Normal"
$ws3.Cells.Item(141, 2).Value = 78
